$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 4-9 (years 2015-2020), columns B-K
$data = @{
    4 = @(17.149638235, 55.1240248708, 14.7144266632, 9.446481546699999, 112.6077535028, 28.9298946537, 104.3657877853, 10.4782592771, 11.0706473679, 1.3152458043)
    5 = @(16.542333, 58.096552, 14.822911, 10.171681, 111.941992, 28.991603, 107.460693, 10.665439, 10.97453, 1.308814)
    6 = @(16.4806991811, 59.8957601621, 14.8140482373, 9.7371215337, 109.7459189229, 29.212351815, 106.6663630277, 10.9109388128, 10.6731617358, 1.3007790158)
    7 = @(16.528363747, 61.9830228564, 14.3189338101, 9.826628323, 109.9724993605, 31.1889245595, 103.0768810615, 10.7567167786, 9.4240010776, 1.2753433894)
    8 = @(16.6517618266, 66.80412564940001, 16.7102769139, 11.4175935958, 110.6098078824, 28.6854542418, 105.8253524596, 11.4931479388, 9.1970798763, 1.150784055)
    9 = @(17.3494231524248, 65.8703711874051, 16.6334152018643, 12.9897280053113, 120.177422164818, 27.373432498469, 109.779031089268, 13.5123244227791, 9.90937496563544, 1.16803442983011)
}

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# New row 10 (year 2021)
$ws.Range("A10").Value = "2021年"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$row10 = @(18.2, 67.7, 16.7, 12.3, 124.8, 34.4, 112, 13.4, 10.1, 1.1)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])10").Value = $row10[$i]
}
